$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the three survey-response strings that were reworded.
$ws.Range("B5").Value = "I have very often experienced such issues in the homes where I've stayed."
$ws.Range("B7").Value = "In the homes where I've lived, I have experienced none of the issues mentioned in the story."
$ws.Range("B16").Value = "I have very often experienced temperature problems in the homes where I have been."

# Match the author's final selection / scroll position.
$ws.Activate()
$ws.Range("B16").Select()
$excel.ActiveWindow.ScrollColumn = 2
$excel.ActiveWindow.ScrollRow = 1
